$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that look numeric (e.g. "1.00",
# "28.332.05"), but the sheet stores them as text. Force the cells we
# are about to rewrite to a text format first so Excel does not silently
# convert strings like "1.00" into the number 1 or "0.0632" into
# scientific notation.
$dCells = @("D2","D3","D5","D6","D8","D11","D12","D13","D14","D16","D17","D18","D19","D20","D23","D24","D25","D26","D27","D30","D33","D36","D38","D39","D41","D43","D44","D47","D48","D50")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '28.332.05'
$ws.Range("E2").Value = '  +4.03%  '
$ws.Range("D3").Value = '1.713.88'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '218.72'
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '23.76'
$ws.Range("E8").Value = '  +2.45%  '
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("D11").Value = '0.0892'
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '1.957.44'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '1.714.64'
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("D14").Value = '4.20'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '67.39'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '28.268.72'
$ws.Range("E17").Value = '  +3.83%  '
$ws.Range("D18").Value = '247.01'
$ws.Range("E18").Value = '  +4.10%  '
$ws.Range("D19").Value = '0.0₃0746'
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = '7.76'
$ws.Range("E20").Value = '  -4.32%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '9.61'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").Value = '147.84'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '7.39'
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").Value = '16.46'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = '0.0511'
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("E31").Value = '  +2.79%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '1.474.27'
$ws.Range("E33").Value = '  -5.18%  '
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("D36").Value = '0.973'
$ws.Range("E36").Value = '  +2.42%  '
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").Value = '0.598'
$ws.Range("E38").Value = '  -1.35%  '
$ws.Range("D39").Value = '0.0175'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").Value = '69.35'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").Value = '5.65'
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("D44").Value = '1.860.30'
$ws.Range("E44").Value = '  +1.18%  '
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").Value = '1.72'
$ws.Range("E47").Value = '  +5.40%  '
$ws.Range("D48").Value = '90.03'
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("E49").Value = '  -2.83%  '
$ws.Range("D50").Value = '8.07'
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("E51").Value = '  -1.74%  '
